$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 14 (ind_0113): B14 -> 0, C14 unchanged
$ws.Range("B14").Value = 0

# row 15 (ind_0114): B15 -> -1, C15 -> -5
$ws.Range("B15").Value = -1
$ws.Range("C15").Value = -5

# row 16 (ind_0115): B16 -> -1, C16 -> -5
$ws.Range("B16").Value = -1
$ws.Range("C16").Value = -5
